$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "68-10=58",
    "2+57=59",
    "45+12=57",
    "35+3=38",
    "31-17=14",
    "35+10=45",
    "18-10=8",
    "5+64=69",
    "1+79=80",
    "90-6=84",
    "48+16=64",
    "47-43=4",
    "20+4=24",
    "80-5=75",
    "66-22=44",
    "50+31=81",
    "69-63=6",
    "46+52=98",
    "67+30=97",
    "84-77=7",
    "70-26=44",
    "20+39=59",
    "35+15=50",
    "78-75=3",
    "79-20=59",
    "90-52=38",
    "65-26=39",
    "80-26=54",
    "66-30=36",
    "34+9=43",
    "19+20=39",
    "80-65=15",
    "37-15=22",
    "38-32=6",
    "71+11=82",
    "63-13=50",
    "52+42=94",
    "33-32=1",
    "37-33=4",
    "39+14=53",
    "79-74=5",
    "59-16=43",
    "42+16=58",
    "46-0=46",
    "89-38=51",
    "6+55=61",
    "10+73=83",
    "75-52=23",
    "71-32=39",
    "93-57=36",
    "39+14=53",
    "12+51=63",
    "46+16=62",
    "33-1=32",
    "5+37=42",
    "66-50=16",
    "20-16=4",
    "59-34=25",
    "73-61=12",
    "28-4=24",
    "60-30=30",
    "73-33=40",
    "3+7=10",
    "91-76=15",
    "92+7=99",
    "64-58=6",
    "91-58=33",
    "88-71=17",
    "96-84=12",
    "8-2=6",
    "29-16=13",
    "43-24=19",
    "4+62=66",
    "18+39=57",
    "79-37=42",
    "95-20=75",
    "88-59=29",
    "13+71=84",
    "74+11=85",
    "31+64=95",
    "31+43=74",
    "15+56=71",
    "45-38=7",
    "76-45=31",
    "30-23=7",
    "13+27=40",
    "27+37=64",
    "9+74=83",
    "78+3=81",
    "55+43=98",
    "41+7=48",
    "54+40=94",
    "29+4=33",
    "76-30=46",
    "78-16=62",
    "16+82=98",
    "81-39=42",
    "83+14=97",
    "33-13=20",
    "63+14=77"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output ("Updated cells: " + $idx)
